# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) for the detail rows (16-38) listed the
# account-statement periods from the most recent (1807) down to the
# oldest (1608). This update reverses that list so it now reads from
# the oldest period (1608) up through the most recent (1807), reflecting
# the newly added part 1 of account-statement periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1608","1609","1610","1611","1612","1701","1702","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807")

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}
